$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new player "Thomas" picks to the list (rows 30-33)
# Enter names/players first (columns A & B), matching the order the
# author typed the rows, then the fbref links (column C) afterwards.
$ws.Range("A30").Value = "Thomas"
$ws.Range("B30").Value = "Pagis"

$ws.Range("A31").Value = "Thomas"
$ws.Range("B31").Value = "Egan Riley"

$ws.Range("A32").Value = "Thomas"
$ws.Range("B32").Value = "Tessman"

$ws.Range("A33").Value = "Thomas"
$ws.Range("B33").Value = "Kechta"

$ws.Range("C30").Value = "https://fbref.com/en/players/7d827b4f/Pablo-Pagis"
$ws.Range("C31").Value = "https://fbref.com/en/players/d313e8ff/CJ-Egan-Riley"
$ws.Range("C32").Value = "https://fbref.com/en/players/ac277993/Tanner-Tessmann"
$ws.Range("C33").Value = "https://fbref.com/en/players/abefc3af/Yassine-Kechta"

# Update the selection to match the last edited cell
$ws.Range("C33").Select()

$wb.Save()
